# "Update countries & provincias Spain"
#
# The workbook is a COVID-19 country leaderboard ("Pais" sheet), sorted by
# total cases descending. This refresh:
#   1) bumps the "last updated" timestamp in A1,
#   2) re-ranks four small clusters of countries whose totals crossed each
#      other since the last refresh (their row keeps its rank/position, but
#      the country name + stats that live at that rank change), and
#   3) applies plain statistic updates to several countries whose rank did
#      not change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Refresh timestamp (row 1) ---------------------------------------
$ws.Range("A1").Value = 'Datos actualizados a 25 de Mayo de 2020 a las 15:35'

# --- 3) Plain stat refreshes (rank/country unchanged) -------------------
# Estados Unidos
$ws.Range("B4").Value = 1688709
$ws.Range("C4").Value = 2273
$ws.Range("E4").Value = 1137616

# India
$ws.Range("B13").Value = 140215
$ws.Range("C13").Value = 1679
$ws.Range("D13").Value = 58216
$ws.Range("E13").Value = 77958

# Panama
$ws.Range("B56").Value = 8360
$ws.Range("C56").Value = 8
$ws.Range("E56").Value = 398

# --- 2a) Yibuti overtakes Bulgaria / Bosnia y Herzegovina / Costa de Marfil (rows 82-85) ---
$ws.Range("A82").Value = 'Republica de Yibuti'
$ws.Range("B82").Value = 2468
$ws.Range("C82").Value = 198
$ws.Range("D82").Value = 1079
$ws.Range("E82").Value = 1375
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 14

$ws.Range("A83").Value = 'Bulgaria'
$ws.Range("B83").Value = 2433
$ws.Range("C83").Value = 6
$ws.Range("D83").Value = 862
$ws.Range("E83").Value = 1441
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 130

$ws.Range("A84").Value = 'Bosnia y Herzegovina'
$ws.Range("B84").Value = 2406
$ws.Range("C84").Value = 5
$ws.Range("D84").Value = 1696
$ws.Range("E84").Value = 564
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 146

$ws.Range("A85").Value = 'Costa de Marfil'
$ws.Range("B85").Value = 2376
$ws.Range("D85").Value = 1219
$ws.Range("E85").Value = 1127
$ws.Range("H85").Value = 30

# --- 3) Cuba (row 90) plain stat refresh ---------------------------------
$ws.Range("B90").Value = 1947
$ws.Range("C90").Value = 6
$ws.Range("D90").Value = 1704
$ws.Range("E90").Value = 161

# --- 2b) Mayotte overtakes Somalia (rows 95-96) --------------------------
$ws.Range("A95").Value = 'Mayotte'
$ws.Range("B95").Value = 1609
$ws.Range("C95").Value = 22
$ws.Range("D95").Value = 894
$ws.Range("E95").Value = 695
$ws.Range("H95").Value = 20

$ws.Range("A96").Value = 'Somalia'
$ws.Range("B96").Value = 1594
$ws.Range("D96").Value = 204
$ws.Range("E96").Value = 1329
$ws.Range("H96").Value = 61

# --- 3) Plain stat refreshes (rank/country unchanged) --------------------
$ws.Range("B103").Value = 1166
$ws.Range("C103").Value = 25
$ws.Range("E103").Value = 461

$ws.Range("B119").Value = 832
$ws.Range("C119").Value = 18
$ws.Range("E119").Value = 108

# --- 2c) Santo Tome y Principe overtakes Nicaragua / Liberia (rows 148-150) ---
$ws.Range("A148").Value = 'Santo Tome y Principe'
$ws.Range("B148").Value = 299
$ws.Range("C148").Value = 48
$ws.Range("D148").Value = 4
$ws.Range("E148").Value = 284
$ws.Range("G148").Value = 3
$ws.Range("H148").Value = 11

$ws.Range("A149").Value = 'Nicaragua'
$ws.Range("B149").Value = 279
$ws.Range("D149").Value = 199
$ws.Range("E149").Value = 63
$ws.Range("H149").Value = 17

$ws.Range("A150").Value = 'Liberia'
$ws.Range("B150").Value = 265
$ws.Range("D150").Value = 141
$ws.Range("E150").Value = 98
$ws.Range("H150").Value = 26

# --- 2d) Malaui overtakes Aruba / Bahamas / Monaco / Barbados / Comoras (rows 170-175) ---
$ws.Range("A170").Value = 'Malaui'
$ws.Range("C170").Value = 18
$ws.Range("D170").Value = 33
$ws.Range("E170").Value = 64
$ws.Range("H170").Value = 4

$ws.Range("A171").Value = 'Aruba'
$ws.Range("B171").Value = 101
$ws.Range("D171").Value = 95
$ws.Range("E171").Value = 3
$ws.Range("H171").Value = 3

$ws.Range("A172").Value = 'Bahamas'
$ws.Range("B172").Value = 100
$ws.Range("D172").Value = 46
$ws.Range("E172").Value = 43
$ws.Range("H172").Value = 11

$ws.Range("A173").Value = 'Monaco'
$ws.Range("B173").Value = 98
$ws.Range("D173").Value = 90
$ws.Range("E173").Value = 4
$ws.Range("H173").Value = 4

$ws.Range("A174").Value = 'Barbados'
$ws.Range("B174").Value = 92
$ws.Range("D174").Value = 70
$ws.Range("E174").Value = 15
$ws.Range("H174").Value = 7

$ws.Range("A175").Value = 'Comoras'
$ws.Range("B175").Value = 87
$ws.Range("D175").Value = 21
$ws.Range("E175").Value = 65
$ws.Range("H175").Value = 1

Write-Host "Update complete"
